# Update Name of Algo
# Apply updated imputed values to result_data_KNN.xlsx (Sheet1)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "C3"   = -13.543
    "B12"  = 5.043
    "C14"  = -12.21
    "C26"  = -13.04
    "B27"  = 5.542000000000001
    "C31"  = -12.629
    "B32"  = 6.776000000000001
    "C35"  = -12.311
    "B36"  = 8.673999999999999
    "C37"  = -13.387
    "B38"  = 5.540000000000001
    "C45"  = -12.776
    "B46"  = 6.366999999999999
    "C52"  = -10.937
    "B54"  = 5.006
    "B55"  = 4.742
    "B56"  = 4.49
    "C57"  = -13.697
    "B67"  = 5.224
    "B69"  = 5.202
    "B72"  = 5.351000000000001
    "C81"  = -13.319
    "B83"  = 5.315
    "C83"  = -14.127
    "B86"  = 4.94
    "B91"  = 6.109
    "B93"  = 5.616
    "B99"  = 5.718000000000001
    "C100" = -13.144
    "C102" = -12.916
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$wb.Save()
